# The underlying TPM values were recomputed; after recalculation, rows whose
# "Target cluster" is "ECs" are no longer part of the output and are dropped,
# while the remaining sending/target cluster pairs keep the refreshed metrics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the trailing rows so the table shrinks from 9 to 6 data rows (A1:T10 -> A1:T7).
$ws.Rows("8:10").Delete()

# Flattened A2:T7 replacement values (20 columns per data row, 6 data rows).
$newValues = @(
    # row 2
    "ECs","L1cam","Alcam","FAPs",3,1,3.685507,11.056521,0.3585631737883472,0.3585631737883472,3,1,0.5683613333333334,1.705084,0.4361027177196302,0.4361027177196302,2.094699672529333,18.852297052764,0.1563703745632743,0.1563703745632743,
    # row 3
    "ECs","L1cam","Alcam","MuSCs",3,1,3.685507,11.056521,0.3585631737883472,0.3585631737883472,3,1,0.7349126666666667,2.204738,0.5638972822803697,0.5638972822803697,2.708525777388667,24.376731996498,0.2021927992250729,0.2021927992250729,
    # row 4
    "FAPs","L1cam","Alcam","FAPs",1,0.3333333333333333,0.099159,0.297477,0.009647184430711629,0.009647184430711629,3,1,0.5683613333333334,1.705084,0.4361027177196302,0.4361027177196302,0.056358141452,0.507223273068,0.004207163348575845,0.004207163348575845,
    # row 5
    "FAPs","L1cam","Alcam","MuSCs",1,0.3333333333333333,0.099159,0.297477,0.009647184430711629,0.009647184430711629,3,1,0.7349126666666667,2.204738,0.5638972822803697,0.5638972822803697,0.072873205114,0.655858846026,0.005440021082135783,0.005440021082135783,
    # row 6
    "MuSCs","L1cam","Alcam","FAPs",3,1,6.493877,19.481631,0.6317896417809412,0.6317896417809411,3,1,0.5683613333333334,1.705084,0.4361027177196302,0.4361027177196302,3.690868590222667,33.217817312004,0.2755251798077801,0.2755251798077801,
    # row 7
    "MuSCs","L1cam","Alcam","MuSCs",3,1,6.493877,19.481631,0.6317896417809412,0.6317896417809411,3,1,0.7349126666666667,2.204738,0.5638972822803697,0.5638972822803697,4.772432463075334,42.951892167678,0.3562644619731611,0.356264461973161
)

$colCount = 20
$rowCount = $newValues.Count / $colCount
for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $newValues[$r * $colCount + $c]
    }
}
